$d = $word.ActiveDocument

$replacements = @(
    @("84÷8=", "92÷4="),
    @("76÷4=", "77÷5="),
    @("34÷5=", "72÷2="),
    @("71÷7=", "53÷7="),
    @("78÷6=", "61÷2="),
    @("55÷6=", "82÷4="),
    @("44÷5=", "78÷7="),
    @("35÷2=", "18÷5="),
    @("30÷2=", "85÷7="),
    @("94÷2=", "15÷9="),
    @("15÷5=", "77÷9="),
    @("37÷6=", "67÷2="),
    @("50÷5=", "60÷4="),
    @("67÷8=", "67÷2="),
    @("28÷9=", "78÷2="),
    @("18÷2=", "54÷2="),
    @("96÷5=", "53÷2="),
    @("73÷6=", "67÷5="),
    @("60÷2=", "22÷6="),
    @("80÷9=", "74÷8="),
    @("87÷2=", "67÷3="),
    @("70÷4=", "57÷7="),
    @("57÷3=", "16÷4="),
    @("21÷6=", "97÷5="),
    @("52÷6=", "20÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
